$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "B2" "1.00"
Set-TextValue "D2" "1.00"

Set-TextValue "B3" "1.00"
Set-TextValue "D3" "1.00"

Set-TextValue "B4" "3.00"
Set-TextValue "D4" "3.00"

Set-TextValue "B6" "9.00"
Set-TextValue "D6" "9.00"
